# AcumulacionResultados.xlsx - header update
#
# "Se actualiza encabezado de plantilla de AcumuladosResultados."
#
# The "Zona Estratégica" column header (I1) is renamed to "Tipo Sitio".
# The old "Zona Estratégica" shared string becomes unused and is dropped
# on save, while the new "Tipo Sitio" string is appended to the shared
# string table - exactly mirroring the canonical OOXML diff (sharedStrings
# entries shift, I1's cached <v> index moves from 20 to the new 25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "Tipo Sitio"

# Leave the cursor on the header cell that was just edited, and let row 1
# grow to fit the (now slightly different) wrapped header text - matches
# the refreshed view state / row height captured in the saved workbook.
$ws.Range("I1").Select() | Out-Null
$ws.Rows.Item(1).RowHeight = 60
